$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Update Price (D) and Volume(1h) (E) columns for rows with value-only changes
Set-TextValue $ws.Range("D2") "24.279.04"
Set-TextValue $ws.Range("E2") "  +14.85%  "
Set-TextValue $ws.Range("D3") "1.680.68"
Set-TextValue $ws.Range("E3") "  +9.33%  "
Set-TextValue $ws.Range("D4") "1.002"
Set-TextValue $ws.Range("E4") "  -0.88%  "
Set-TextValue $ws.Range("D5") "307.27"
Set-TextValue $ws.Range("E5") "  +8.99%  "
Set-TextValue $ws.Range("D6") "0.9967"
Set-TextValue $ws.Range("E6") "  +3.11%  "
Set-TextValue $ws.Range("D7") "0.3722"
Set-TextValue $ws.Range("E7") "  +2.72%  "
Set-TextValue $ws.Range("D8") "0.3430"
Set-TextValue $ws.Range("E8") "  +8.18%  "
Set-TextValue $ws.Range("D9") "48.24"
Set-TextValue $ws.Range("E9") "  +18.58%  "
Set-TextValue $ws.Range("D10") "1.183"
Set-TextValue $ws.Range("E10") "  +8.18%  "
Set-TextValue $ws.Range("D11") "0.07292"
Set-TextValue $ws.Range("E11") "  +7.04%  "
Set-TextValue $ws.Range("D12") "0.9967"
Set-TextValue $ws.Range("E12") "  -0.83%  "
Set-TextValue $ws.Range("D13") "20.57"
Set-TextValue $ws.Range("E13") "  +9.94%  "
Set-TextValue $ws.Range("D14") "6.090"
Set-TextValue $ws.Range("E14") "  +7.56%  "
Set-TextValue $ws.Range("D15") "6.746"
Set-TextValue $ws.Range("E15") "  +6.27%  "
Set-TextValue $ws.Range("D16") "1.677.68"
Set-TextValue $ws.Range("E16") "  +9.80%  "
Set-TextValue $ws.Range("D17") "0.00001107"
Set-TextValue $ws.Range("E17") "  +6.00%  "
Set-TextValue $ws.Range("D18") "0.9964"
Set-TextValue $ws.Range("E18") "  +3.12%  "
Set-TextValue $ws.Range("D19") "0.06713"
Set-TextValue $ws.Range("E19") "  +10.46%  "
Set-TextValue $ws.Range("D20") "81.48"
Set-TextValue $ws.Range("E20") "  +12.64%  "
Set-TextValue $ws.Range("D21") "16.45"
Set-TextValue $ws.Range("E21") "  +9.77%  "
Set-TextValue $ws.Range("D22") "6.122"
Set-TextValue $ws.Range("E22") "  +7.38%  "
Set-TextValue $ws.Range("D23") "12.03"
Set-TextValue $ws.Range("E23") "  +5.86%  "
Set-TextValue $ws.Range("D24") "24.243.34"
Set-TextValue $ws.Range("E24") "  +14.51%  "
Set-TextValue $ws.Range("D25") "2.400"
Set-TextValue $ws.Range("E25") "  +3.48%  "
Set-TextValue $ws.Range("D28") "152.14"
Set-TextValue $ws.Range("E28") "  +2.56%  "
Set-TextValue $ws.Range("D29") "19.52"
Set-TextValue $ws.Range("E29") "  +10.68%  "
Set-TextValue $ws.Range("D30") "1.856.80"
Set-TextValue $ws.Range("E30") "  +9.52%  "
Set-TextValue $ws.Range("D31") "127.00"
Set-TextValue $ws.Range("E31") "  +7.18%  "
Set-TextValue $ws.Range("D32") "6.418"
Set-TextValue $ws.Range("E32") "  +24.00%  "
Set-TextValue $ws.Range("D33") "4.024"
Set-TextValue $ws.Range("E33") "  -0.05%  "
Set-TextValue $ws.Range("D34") "0.9870"
Set-TextValue $ws.Range("E34") "  +16.07%  "
Set-TextValue $ws.Range("D35") "1.740"
Set-TextValue $ws.Range("E35") "  +15.55%  "
Set-TextValue $ws.Range("D36") "0.08448"
Set-TextValue $ws.Range("E36") "  +5.69%  "
Set-TextValue $ws.Range("D38") "5.368"
Set-TextValue $ws.Range("E38") "  +8.40%  "
Set-TextValue $ws.Range("D39") "0.06410"
Set-TextValue $ws.Range("E39") "  +9.09%  "
Set-TextValue $ws.Range("D40") "8.846"
Set-TextValue $ws.Range("E40") "  +14.60%  "
Set-TextValue $ws.Range("D41") "1.290"
Set-TextValue $ws.Range("E41") "  +6.59%  "
Set-TextValue $ws.Range("D42") "0.02334"
Set-TextValue $ws.Range("E42") "  +11.16%  "
Set-TextValue $ws.Range("D43") "0.2110"
Set-TextValue $ws.Range("E43") "  +10.34%  "
Set-TextValue $ws.Range("D44") "0.6142"
Set-TextValue $ws.Range("E44") "  +13.01%  "
Set-TextValue $ws.Range("D45") "0.9961"
Set-TextValue $ws.Range("E45") "  +3.05%  "
Set-TextValue $ws.Range("D46") "3.796"
Set-TextValue $ws.Range("E46") "  +6.35%  "
Set-TextValue $ws.Range("D47") "13.19"
Set-TextValue $ws.Range("E47") "  +5.36%  "
Set-TextValue $ws.Range("D48") "0.5953"
Set-TextValue $ws.Range("E48") "  +9.55%  "
Set-TextValue $ws.Range("D49") "128.07"
Set-TextValue $ws.Range("E49") "  +5.57%  "
Set-TextValue $ws.Range("D50") "2.018"
Set-TextValue $ws.Range("E50") "  +8.03%  "
Set-TextValue $ws.Range("D51") "0.07158"
Set-TextValue $ws.Range("E51") "  +8.90%  "

# Row 37: only Volume(1h) changes
Set-TextValue $ws.Range("E37") "  +17.21%  "

# Rows 26 and 27 swap content: LidoDAOToken <-> LEO
Set-TextValue $ws.Range("B26") "LEO"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D26") "3.360"
Set-TextValue $ws.Range("E26") "  -9.08%  "

Set-TextValue $ws.Range("B27") "LidoDAOToken"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D27") "2.669"
Set-TextValue $ws.Range("E27") "  +20.50%  "
